$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary = $wb.Worksheets.Item("Summary")

$wsForecast.Range("D11").Value = 63
$wsForecast.Range("D12").Value = 73
$wsForecast.Range("D13").Value = 73
$wsForecast.Range("D14").Value = 67
$wsForecast.Range("D15").Value = 66
$wsForecast.Range("D16").Value = 79

# B9 on the Summary sheet holds the "Total Forecast (16 Weeks)" total as text,
# mirroring the original inline-string cell, so force Text formatting before
# writing the new value and then clear the number-format override again.
$wsSummary.Range("B9").NumberFormat = "@"
$wsSummary.Range("B9").Value = "1376"
$wsSummary.Range("B9").ClearFormats()
